# Insert a new record (row) for Berenjena at Femacal de La Calera just
# after the existing row 310, pushing the former rows 311..431 down to
# 312..432 (dimension grows from A1:R431 to A1:R432).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("311:311").Insert()

$ws.Range("A311").Value = 3
$ws.Range("B311").Value = "Femacal de La Calera"
$ws.Range("C311").Value = "Coquimbo"
$ws.Range("D311").Value = 45027
$ws.Range("E311").Value = 5
$ws.Range("F311").Value = 100112001
$ws.Range("G311").Value = "Berenjena"
$ws.Range("H311").Value = "Sin especificar"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 85
$ws.Range("K311").Value = 7500
$ws.Range("L311").Value = 8000
$ws.Range("M311").Value = 7735
$ws.Range("N311").Value = "$/caja 60 unidades"
$ws.Range("O311").Value = "Región de Arica y Parinacota"
$ws.Range("P311").Value = 129
$ws.Range("Q311").Value = 60
$ws.Range("R311").Value = "Hortaliza"
